$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "sensactIO4" column (table header / G1) to "sensactHsNano"
$ws.Range("G1").Value = "sensactHsNano"

# The sensactHsNano board's pin assignment for this column hasn't been
# filled in yet (moving from ESP32DeviceController's sensactIO4 PCB) -
# clear out all the old per-pin notes in column G, keeping row 17's
# "Gut" highlight style so it stays ready for new data.
$ws.Range("G2:G35").ClearContents()

# Reflect the selection left behind after clearing the column
$ws.Range("G2:G35").Select()
